# "actualizacion de formulas para integrales"
# Adds the MOD/IF helper formula in E8 and two blank-padding label cells
# (H15, I15) used as spacers, then leaves the selection on I9 to match
# the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New boolean-valued helper formula (row 8) - mirrors the other IF(...)
# helper formulas already present in row 4, but with no FALSE branch so a
# failed MOD test yields boolean FALSE.
$ws.Range("E8").Formula = "=IF(MOD(C3,2)=0,(F3/3)*C5+(4*D5)+(4*E5)+(4*F5)+(4*G5)+H5)"

# New row 15: a single-space and a triple-space text cell.
$ws.Range("H15").Value = " "
$ws.Range("I15").Value = "   "

# Leave the selection where the author left it.
[void]$ws.Range("I9").Select()
